$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @{
    2 = @("28.417.91", "+3.47%")
    3 = @("1.870.97", "+2.12%")
    4 = @($null, "-0.26%")
    5 = @("339.03", "+2.53%")
    6 = @($null, "-0.27%")
    7 = @("0.4703", "+2.33%")
    8 = @("0.3973", "+4.13%")
    9 = @("47.77", "+2.24%")
    10 = @("0.08041", "+1.62%")
    11 = @($null, "+3.33%")
    12 = @("22.03", "+4.77%")
    13 = @("6.052", "+2.77%")
    14 = @("1.872.15", "+2.17%")
    15 = @("7.265", "+3.29%")
    16 = @("91.21", "+3.66%")
    17 = @("1.001", "-0.37%")
    18 = @("0.00001040", "+1.12%")
    19 = @("0.06628", "+0.07%")
    20 = @("17.61", "+3.82%")
    21 = @($null, "-0.28%")
    22 = @("28.424.18", "+3.50%")
    23 = @("5.478", "+2.52%")
    24 = @($null, "+2.35%")
    25 = @("2.253", $null)
    26 = @("2.092.12", "+1.93%")
    27 = @("160.66", "+2.03%")
    28 = @("19.77", "+2.06%")
    29 = @("2.127", "+3.19%")
    30 = @("5.512", "+3.91%")
    31 = @("120.41", "+1.23%")
    32 = @("0.9750", "+2.10%")
    33 = @("0.09514", "+2.53%")
    34 = @("3.592", "+0.18%")
    35 = @("1.377", "+4.88%")
    36 = @("5.353", "+2.11%")
    37 = @("0.06108", "+3.08%")
    38 = @("0.02256", "+2.88%")
    39 = @("8.368", "+3.92%")
    40 = @("1.179", "+1.26%")
    41 = @("0.5962", "+2.98%")
    42 = @($null, "-0.30%")
    43 = @("0.1882", "+2.41%")
    44 = @("10.35", "+3.21%")
    45 = @($null, "+2.94%")
    46 = @("0.5599", "+2.18%")
    47 = @("12.19", "+1.82%")
    48 = @("1.958", "+5.04%")
    49 = @("0.06937", "+4.45%")
    50 = @("2.065", "+16.28%")
    51 = @("111.73", $null)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $dVal = $vals[0]
    $eVal = $vals[1]
    if ($dVal -ne $null) {
        $ws.Cells.Item($row, 4).Value = $dVal
    }
    if ($eVal -ne $null) {
        $ws.Cells.Item($row, 5).Value = "  " + $eVal + "  "
    }
}
